# Rename column headers in row 1 from Chinese to English during preprocessing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Field ID"
$ws.Range("B1").Value = "Crop ID"
$ws.Range("C1").Value = "Crop Name"
$ws.Range("D1").Value = "Crop Type"
$ws.Range("E1").Value = "Planting Area"
$ws.Range("F1").Value = "Season"
$ws.Range("G1").Value = "Field Type"
$ws.Range("H1").Value = "Field Area"
$ws.Range("I1").Value = "Yield"
$ws.Range("J1").Value = "Cost"
$ws.Range("K1").Value = "Price"
